$d = $word.ActiveDocument

# 1. The original document opens with a "Trash pickup" / "Municipal trash
#    collection..." pair of paragraphs that carries the _GoBack bookmark,
#    immediately followed by the "(10 points)" user story. The edit moves
#    that lead pair one slot further down, past the bookmark, so the
#    "(10 points)" story paragraph becomes the new first paragraph in the
#    body and also the new home of the _GoBack bookmark. Cut/Paste (rather
#    than deleting + retyping the text) relocates the two paragraphs while
#    keeping their original, unadorned run/paragraph markup intact.
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$leadRange = $d.Range($p1.Range.Start, $p2.Range.End)
$leadRange.Cut()

$top = $d.Range(0, 0)
$top.Paste()

# 2. Relocate the _GoBack bookmark so it sits immediately before the
#    "(10 points)" run instead of on the original "Trash pickup" paragraph.
#    Word only keeps a single _GoBack bookmark, so adding the new one
#    implicitly removes the old one (which Cut() above already dropped).
$findRng = $d.Content
$findRng.Find.Execute("(10 points) As a new visitor", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($findRng.Find.Found) {
    $bmStart = $findRng.Start
    $bmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# 3. Highlight (yellow) the three user-story paragraphs that call out the
#    sign-up flow, the pickup-day change, and the vacation hold.
$targets = @(
  "(10 points) As a new visitor to the site, I want to be able to enter my info and sign up so I can start receiving regular trash pickups.",
  "(15 points) As a registered user, I want to be able to change my pickup day for one week or for every week to fit my schedule.",
  "(20 points) As a registered user, I want to be able to specify a period of time that I don't want pickup so I can go on vacation."
)

foreach ($t in $targets) {
    $rng = $d.Content
    $rng.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($rng.Find.Found) {
        $rng.HighlightColorIndex = 7
    }
}
